# Shrink the login test-data sheet from a 4-column/4-row grid down to a
# simple 2x2 username/password table, and replace the old sample
# "firstname"/"lastname" header values with actual admin credentials
# ("Admin" / "admin123") used by the test suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused C and D columns (firstname/lastname + sample data)
# and the extra sample rows 3 and 4, leaving a tight A1:B2 range.
$ws.Range("C:D").Delete()
$ws.Range("3:4").Delete()

# Replace the row-2 values with the admin credentials.
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

# Match the saved selection state (cell below the data table).
[void]$ws.Range("B3").Select()
